$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data refresh.
# D-column values are written via an apostrophe-prefixed Formula (then the
# original cell style is restored) so numeric-looking text such as "49.50"
# or "3.800" is preserved verbatim as text instead of being coerced into a
# number (which would silently drop trailing zeros / introduce FP noise).
$style = $ws.Range("D2").Style
$ws.Range("D2").Formula = "'22.376.67"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  +0.00%  "
$style = $ws.Range("D3").Style
$ws.Range("D3").Formula = "'1.566.54"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.11%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").Formula = "'1.002"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +0.03%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").Formula = "'291.24"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +0.61%  "
$style = $ws.Range("D7").Style
$ws.Range("D7").Formula = "'0.3761"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  +2.19%  "
$style = $ws.Range("D8").Style
$ws.Range("D8").Formula = "'49.50"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("E9").Value = "  +0.83%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").Formula = "'0.07598"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -0.33%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").Formula = "'1.139"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -2.08%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").Formula = "'1.003"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +0.06%  "
$style = $ws.Range("D13").Style
$ws.Range("D13").Formula = "'21.04"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  -1.21%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").Formula = "'5.982"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -1.29%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").Formula = "'6.948"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +0.56%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").Formula = "'1.570.28"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +0.26%  "
$style = $ws.Range("D17").Style
$ws.Range("D17").Formula = "'0.00001131"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -0.22%  "
$style = $ws.Range("D18").Style
$ws.Range("D18").Formula = "'89.94"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  +0.45%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").Formula = "'0.06735"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("E20").Value = "  +0.06%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").Formula = "'16.58"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("E22").Value = "  -0.67%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").Formula = "'11.93"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -0.31%  "
$style = $ws.Range("D24").Style
$ws.Range("D24").Formula = "'22.369.58"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  +0.53%  "
$style = $ws.Range("D26").Style
$ws.Range("D26").Formula = "'2.689"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -7.85%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").Formula = "'20.06"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +0.30%  "
$style = $ws.Range("D28").Style
$ws.Range("D28").Formula = "'147.29"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  +0.95%  "
$style = $ws.Range("D29").Style
$ws.Range("D29").Formula = "'5.024"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("E30").Value = "  +0.80%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").Formula = "'1.740.77"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -0.28%  "
$style = $ws.Range("D32").Style
$ws.Range("D32").Formula = "'2.018"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +0.47%  "
$style = $ws.Range("D33").Style
$ws.Range("D33").Formula = "'0.9998"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -3.59%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").Formula = "'6.072"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -3.02%  "
$style = $ws.Range("D35").Style
$ws.Range("D35").Formula = "'10.10"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -1.24%  "
$style = $ws.Range("D36").Style
$ws.Range("D36").Formula = "'0.08492"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +0.16%  "
$style = $ws.Range("D37").Style
$ws.Range("D37").Formula = "'0.02512"
$ws.Range("D37").Style = $style
$style = $ws.Range("D38").Style
$ws.Range("D38").Formula = "'1.374"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +6.97%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").Formula = "'0.2298"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -0.99%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").Formula = "'0.06512"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -0.08%  "
$style = $ws.Range("D41").Style
$ws.Range("D41").Formula = "'5.388"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -2.60%  "
$style = $ws.Range("D42").Style
$ws.Range("D42").Formula = "'11.34"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -2.95%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").Formula = "'0.6321"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("E44").Value = "  +0.19%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").Formula = "'13.98"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -2.61%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").Formula = "'3.800"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +1.45%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").Formula = "'0.5926"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -0.86%  "
$style = $ws.Range("D48").Style
$ws.Range("D48").Formula = "'2.074"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -1.43%  "
$style = $ws.Range("D49").Style
$ws.Range("D49").Formula = "'1.273"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  +1.07%  "
$style = $ws.Range("D50").Style
$ws.Range("D50").Formula = "'124.47"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +0.31%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").Formula = "'0.07311"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +0.47%  "
